$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price (column D) updates ---
# Values are stored as text in the sheet; a leading apostrophe forces Excel
# to keep them as text instead of auto-converting to a number.
$ws.Cells.Item(2, 4).Value  = "'244.38"
$ws.Cells.Item(4, 4).Value  = "'5.396"
$ws.Cells.Item(5, 4).Value  = "'0.05987"
$ws.Cells.Item(7, 4).Value  = "'0.8145"
$ws.Cells.Item(8, 4).Value  = "'0.9551"
$ws.Cells.Item(9, 4).Value  = "'0.1425"
$ws.Cells.Item(10, 4).Value = "'0.07442"
$ws.Cells.Item(11, 4).Value = "'0.03282"
$ws.Cells.Item(12, 4).Value = "'0.03054"
$ws.Cells.Item(13, 4).Value = "'0.09411"
$ws.Cells.Item(14, 4).Value = "'4.002"
$ws.Cells.Item(15, 4).Value = "'0.001587"
$ws.Cells.Item(16, 4).Value = "'0.04813"
$ws.Cells.Item(18, 4).Value = "'0.005454"
$ws.Cells.Item(20, 4).Value = "'0.0009883"
$ws.Cells.Item(22, 4).Value = "'3.679"
$ws.Cells.Item(23, 4).Value = "'6.434"
$ws.Cells.Item(24, 4).Value = "'2.189"
$ws.Cells.Item(40, 4).Value = "'0.03999"

# --- Rows 41-43: coin list reshuffled (KickToken/BKEXToken/CEJI rotate) ---
# Row 41 becomes BKEXToken
$ws.Cells.Item(41, 2).Value = "BKEXToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Cells.Item(41, 4).Value = "'0.1074"
$ws.Cells.Item(41, 5).Value = "40BKEXTokenBKK"

# Row 42 becomes CEJI
$ws.Cells.Item(42, 2).Value = "CEJI"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Cells.Item(42, 4).Value = "'0.002721"
$ws.Cells.Item(42, 5).Value = "41CEJICEJI"

# Row 43 becomes KickToken
$ws.Cells.Item(43, 2).Value = "KickToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Cells.Item(43, 4).Value = "'0.006563"
$ws.Cells.Item(43, 5).Value = "42KickTokenKICK"

# --- remaining simple updates ---
$ws.Cells.Item(44, 4).Value = "'0.005712"
$ws.Cells.Item(45, 4).Value = "'0.00005139"
$ws.Cells.Item(47, 5).Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Cells.Item(48, 4).Value = "'0.006306"
$ws.Cells.Item(49, 4).Value = "'0.00002101"
